$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph.
$f = $d.Content
$found = $f.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $verParagraph = $f.Paragraphs(1)

    # The next paragraph is the "(c) 2020 ..." footer line.
    $copyrightParagraph = $verParagraph.Next()

    # The paragraph after that is the (now redundant) blank spacer paragraph
    # that used to sit between the footer text and the trailing page break.
    $blankParagraph = $copyrightParagraph.Next()

    # Delete from the start of "Ver no Jupiter..." through the end of the
    # blank spacer paragraph (its paragraph mark included), which removes
    # both text paragraphs and collapses the extra blank line in one go,
    # leaving a single blank paragraph before the page-break paragraph.
    $rng = $d.Range($verParagraph.Range.Start, $blankParagraph.Range.End)
    $rng.Delete()
}
